$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G30").Value = 0.15
$ws.Range("I30").Value = 3
$ws.Range("K30").Formula = "=G30*I30"
$ws.Range("A30").Value = "Resistor for Relays (Non Connected) (0603) (700ohm)"
$ws.Range("M30").Value = "https://www.digikey.ca/en/products/detail/stackpole-electronics-inc/RMCF0603FT732R/1714224"

Write-Host "K30 value2:"
Write-Host $ws.Range("K30").Value2
